$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change) per upstream data refresh.

$ws.Range("D2").Value = '26.767.49'
$ws.Range("E2").Value = '  -4.23%  '

$ws.Range("D3").Value = '1.718.50'
$ws.Range("E3").Value = '  -2.75%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.05'
$ws.Range("E5").Value = '  -5.92%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4859'
$ws.Range("E7").Value = '  +4.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3475'
$ws.Range("E8").Value = '  -1.26%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.70'
$ws.Range("E9").Value = '  -1.84%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07207'
$ws.Range("E10").Value = '  -2.43%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.044'
$ws.Range("E11").Value = '  -3.65%  '

$ws.Range("E12").Value = '  +0.26%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.72'
$ws.Range("E13").Value = '  -4.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.841'
$ws.Range("E14").Value = '  -2.76%  '

$ws.Range("D15").Value = '1.724.12'
$ws.Range("E15").Value = '  -2.43%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.788'
$ws.Range("E16").Value = '  -5.52%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.07'
$ws.Range("E17").Value = '  -6.67%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001032'
$ws.Range("E18").Value = '  -2.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06401'
$ws.Range("E19").Value = '  -0.37%  '

$ws.Range("E20").Value = '  +0.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.47'
$ws.Range("E21").Value = '  -2.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.689'
$ws.Range("E22").Value = '  -1.73%  '

$ws.Range("D23").Value = '26.841.67'
$ws.Range("E23").Value = '  -4.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.87'
$ws.Range("E24").Value = '  -2.24%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.051'
$ws.Range("E25").Value = '  -4.84%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.32'
$ws.Range("E26").Value = '  -5.83%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.77'
$ws.Range("E27").Value = '  -1.24%  '

$ws.Range("D28").Value = '1.921.08'
$ws.Range("E28").Value = '  -2.45%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.050'
$ws.Range("E29").Value = '  -6.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.80'
$ws.Range("E30").Value = '  -2.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.025'
$ws.Range("E31").Value = '  -4.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09278'
$ws.Range("E32").Value = '  -0.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.561'
$ws.Range("E33").Value = '  -2.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.335'
$ws.Range("E34").Value = '  -3.83%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05867'
$ws.Range("E35").Value = '  -3.87%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02170'
$ws.Range("E36").Value = '  -4.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.426'
$ws.Range("E37").Value = '  -1.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.88'
$ws.Range("E38").Value = '  -6.95%  '

$ws.Range("E39").Value = '  +0.15%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.1977'
$ws.Range("E40").Value = '  -4.40%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.713'
$ws.Range("E41").Value = '  -4.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5945'
$ws.Range("E42").Value = '  -3.42%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.115'
$ws.Range("E43").Value = '  -6.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.400'
$ws.Range("E44").Value = '  -4.79%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.77'
$ws.Range("E45").Value = '  -2.64%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.568'
$ws.Range("E46").Value = '  -4.70%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5568'
$ws.Range("E47").Value = '  -3.90%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '118.68'
$ws.Range("E48").Value = '  -4.16%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.824'
$ws.Range("E49").Value = '  -5.57%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06615'
$ws.Range("E50").Value = '  -2.96%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.089'
$ws.Range("E51").Value = '  -3.15%  '
